$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 <- content from original row 3
$ws.Range("A2").Value = 130803083
$ws.Range("B2").Value = 83089
$ws.Range("D2").Value = "NT"
$ws.Range("E2").Value = 1312
$ws.Range("F2").Value = "Gammelgransskål"
$ws.Range("G2").Value = "Pseudographis pinicola"
$ws.Range("H2").Value = "(Nyl.) Rehm"
$ws.Range("Q2").Value = 424802
$ws.Range("R2").Value = 6712148

# Row 3 <- content from original row 4
$ws.Range("A3").Value = 130803092
$ws.Range("B3").Value = 89193
$ws.Range("D3").Value = "NT"
$ws.Range("E3").Value = 510
$ws.Range("F3").Value = "Doftskinn"
$ws.Range("G3").Value = "Cystostereum murrayi"
$ws.Range("H3").Value = "(Berk. & M.A.Curtis.) Pouzar"
$ws.Range("Q3").Value = 424832
$ws.Range("R3").Value = 6712186

# Row 4 <- content from original row 2
$ws.Range("A4").Value = 130803039
$ws.Range("B4").Value = 83223
$ws.Range("D4").Value = "NT"
$ws.Range("E4").Value = 6440
$ws.Range("F4").Value = "Vitgrynig nållav"
$ws.Range("G4").Value = "Chaenotheca subroscida"
$ws.Range("H4").Value = "(Eitner) Zahlbr."
$ws.Range("Q4").Value = 424963
$ws.Range("R4").Value = 6712076

# Row 7 <- content from original row 11
$ws.Range("A7").Value = 130803088
$ws.Range("B7").Value = 83089
$ws.Range("D7").Value = "NT"
$ws.Range("E7").Value = 1312
$ws.Range("F7").Value = "Gammelgransskål"
$ws.Range("G7").Value = "Pseudographis pinicola"
$ws.Range("H7").Value = "(Nyl.) Rehm"
$ws.Range("Q7").Value = 424964
$ws.Range("R7").Value = 6712067

# Row 8 <- content from original row 10
$ws.Range("A8").Value = 130803067
$ws.Range("B8").Value = 78255
$ws.Range("D8").Value = "NT"
$ws.Range("E8").Value = 228579
$ws.Range("F8").Value = "Liten svartspik"
$ws.Range("G8").Value = "Chaenothecopsis nana"
$ws.Range("H8").Value = "Tibell"
$ws.Range("Q8").Value = 424814
$ws.Range("R8").Value = 6712361

# Row 9 <- content from original row 7
$ws.Range("A9").Value = 130803042
$ws.Range("B9").Value = 91771
$ws.Range("D9").Value = "LC"
$ws.Range("E9").Value = 5447
$ws.Range("F9").Value = "Vedticka"
$ws.Range("G9").Value = "Fuscoporia viticola"
$ws.Range("H9").Value = "(Schwein.) Murrill"
$ws.Range("Q9").Value = 424979
$ws.Range("R9").Value = 6712092

# Row 10 <- content from original row 8
$ws.Range("A10").Value = 130803071
$ws.Range("B10").Value = 91181
$ws.Range("D10").Value = "LC"
$ws.Range("E10").Value = 5685
$ws.Range("F10").Value = "Gullgröppa"
$ws.Range("G10").Value = "Pseudomerulius aureus"
$ws.Range("H10").Value = "(Fr.) Jülich"
$ws.Range("Q10").Value = 424873
$ws.Range("R10").Value = 6712126

# Row 11 <- content from original row 9
$ws.Range("A11").Value = 130803064
$ws.Range("B11").Value = 91829
$ws.Range("D11").Value = "NT"
$ws.Range("E11").Value = 5442
$ws.Range("F11").Value = "Tallticka"
$ws.Range("G11").Value = "Porodaedalea pini"
$ws.Range("H11").Value = "(Brot.) Murrill"
$ws.Range("Q11").Value = 424893
$ws.Range("R11").Value = 6712101

# Row 12 <- content from original row 13
$ws.Range("A12").Value = 130803066
$ws.Range("B12").Value = 75221
$ws.Range("D12").Value = "LC"
$ws.Range("E12").Value = 6428
$ws.Range("F12").Value = "Rostfläck"
$ws.Range("G12").Value = "Arthonia vinosa"
$ws.Range("H12").Value = "Leight."
$ws.Range("Q12").Value = 424814
$ws.Range("R12").Value = 6712218

# Row 13 <- content from original row 12
$ws.Range("A13").Value = 130803085
$ws.Range("B13").Value = 83089
$ws.Range("D13").Value = "NT"
$ws.Range("E13").Value = 1312
$ws.Range("F13").Value = "Gammelgransskål"
$ws.Range("G13").Value = "Pseudographis pinicola"
$ws.Range("H13").Value = "(Nyl.) Rehm"
$ws.Range("Q13").Value = 424798
$ws.Range("R13").Value = 6712168

# Row 30 <- content from original row 31
$ws.Range("A30").Value = 130803079
$ws.Range("B30").Value = 79243
$ws.Range("D30").Value = "NT"
$ws.Range("E30").Value = 6425
$ws.Range("F30").Value = "Garnlav"
$ws.Range("G30").Value = "Alectoria sarmentosa"
$ws.Range("H30").Value = "(Ach.) Ach."
$ws.Range("Q30").Value = 424980
$ws.Range("R30").Value = 6712107

# Row 31 <- content from original row 30
$ws.Range("A31").Value = 130803080
$ws.Range("B31").Value = 83089
$ws.Range("D31").Value = "NT"
$ws.Range("E31").Value = 1312
$ws.Range("F31").Value = "Gammelgransskål"
$ws.Range("G31").Value = "Pseudographis pinicola"
$ws.Range("H31").Value = "(Nyl.) Rehm"
$ws.Range("Q31").Value = 424918
$ws.Range("R31").Value = 6712188

# Row 32 <- content from original row 33
$ws.Range("A32").Value = 130803081
$ws.Range("B32").Value = 83089
$ws.Range("D32").Value = "NT"
$ws.Range("E32").Value = 1312
$ws.Range("F32").Value = "Gammelgransskål"
$ws.Range("G32").Value = "Pseudographis pinicola"
$ws.Range("H32").Value = "(Nyl.) Rehm"
$ws.Range("Q32").Value = 424813
$ws.Range("R32").Value = 6712360

# Row 33 <- content from original row 32
$ws.Range("A33").Value = 130803098
$ws.Range("B33").Value = 83215
$ws.Range("D33").Value = "NT"
$ws.Range("E33").Value = 308
$ws.Range("F33").Value = "Brunpudrad nållav"
$ws.Range("G33").Value = "Chaenotheca gracillima"
$ws.Range("H33").Value = "(Vain.) Tibell"
$ws.Range("Q33").Value = 424798
$ws.Range("R33").Value = 6712176

# Row 34 <- content from original row 36
$ws.Range("A34").Value = 130803077
$ws.Range("B34").Value = 79243
$ws.Range("D34").Value = "NT"
$ws.Range("E34").Value = 6425
$ws.Range("F34").Value = "Garnlav"
$ws.Range("G34").Value = "Alectoria sarmentosa"
$ws.Range("H34").Value = "(Ach.) Ach."
$ws.Range("Q34").Value = 424877
$ws.Range("R34").Value = 6712121

# Row 35 <- content from original row 37
$ws.Range("A35").Value = 130803078
$ws.Range("B35").Value = 79243
$ws.Range("D35").Value = "NT"
$ws.Range("E35").Value = 6425
$ws.Range("F35").Value = "Garnlav"
$ws.Range("G35").Value = "Alectoria sarmentosa"
$ws.Range("H35").Value = "(Ach.) Ach."
$ws.Range("Q35").Value = 424951
$ws.Range("R35").Value = 6712137

# Row 36 <- content from original row 35
$ws.Range("A36").Value = 130803082
$ws.Range("B36").Value = 83089
$ws.Range("D36").Value = "NT"
$ws.Range("E36").Value = 1312
$ws.Range("F36").Value = "Gammelgransskål"
$ws.Range("G36").Value = "Pseudographis pinicola"
$ws.Range("H36").Value = "(Nyl.) Rehm"
$ws.Range("Q36").Value = 424785
$ws.Range("R36").Value = 6712249

# Row 37 <- content from original row 34
$ws.Range("A37").Value = 130803095
$ws.Range("B37").Value = 83215
$ws.Range("D37").Value = "NT"
$ws.Range("E37").Value = 308
$ws.Range("F37").Value = "Brunpudrad nållav"
$ws.Range("G37").Value = "Chaenotheca gracillima"
$ws.Range("H37").Value = "(Vain.) Tibell"
$ws.Range("Q37").Value = 424724
$ws.Range("R37").Value = 6712378

# Row 56 <- content from original row 57
$ws.Range("A56").Value = 130803048
$ws.Range("B56").Value = 57884
$ws.Range("D56").Value = "NT"
$ws.Range("E56").Value = 100109
$ws.Range("F56").Value = "Tretåig hackspett"
$ws.Range("G56").Value = "Picoides tridactylus"
$ws.Range("H56").Value = "(Linnaeus, 1758)"
$ws.Range("Q56").Value = 424801
$ws.Range("R56").Value = 6712325

# Row 57 <- content from original row 56
$ws.Range("A57").Value = 130803057
$ws.Range("B57").Value = 57884
$ws.Range("D57").Value = "NT"
$ws.Range("E57").Value = 100109
$ws.Range("F57").Value = "Tretåig hackspett"
$ws.Range("G57").Value = "Picoides tridactylus"
$ws.Range("H57").Value = "(Linnaeus, 1758)"
$ws.Range("Q57").Value = 424935
$ws.Range("R57").Value = 6712079
